$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:E6 and G2:G6 (F column unchanged)
$ws.Range("B2").Value = 0.1554434735375247
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 7.127362168834791

$ws.Range("B3").Value = 0.3464964993005633
$ws.Range("C3").Value = 0.05231270169004087
$ws.Range("D3").Value = 16.98373111632243
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 23.86396839458365

$ws.Range("B4").Value = 1.505614041169197
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 0.1529057820181812
$ws.Range("E4").Value = 6.48142807727062
$ws.Range("G4").Value = 9.793184359356808

$ws.Range("B5").Value = 1.505614041169197
$ws.Range("C5").Value = 9.226618575922256
$ws.Range("D5").Value = 16.98373111632243
$ws.Range("E5").Value = 6.48142807727062
$ws.Range("G5").Value = 34.1973918106845

$ws.Range("B6").Value = 3.182878228561681
$ws.Range("C6").Value = 1.65323645889881
$ws.Range("D6").Value = 0.1529057820181812
$ws.Range("E6").Value = 0.4998867070740569
$ws.Range("G6").Value = 5.488907176552729
